$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 10.32689766666667
$ws.Cells.Item(2, 8).Value = 30.980693
$ws.Cells.Item(2, 9).Value = 0.2044815006034941
$ws.Cells.Item(2, 10).Value = 0.204481500603494
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 16.535604
$ws.Cells.Item(2, 14).Value = 49.606812
$ws.Cells.Item(2, 15).Value = 0.2120453146491552
$ws.Cells.Item(2, 16).Value = 0.2120453146491552
$ws.Cells.Item(2, 17).Value = 170.761490364524
$ws.Cells.Item(2, 18).Value = 1536.853413280716
$ws.Cells.Item(2, 19).Value = 0.04335934413539932
$ws.Cells.Item(2, 20).Value = 0.04335934413539932

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 10.32689766666667
$ws.Cells.Item(3, 8).Value = 30.980693
$ws.Cells.Item(3, 9).Value = 0.2044815006034941
$ws.Cells.Item(3, 10).Value = 0.204481500603494
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 40.62063066666667
$ws.Cells.Item(3, 14).Value = 121.861892
$ws.Cells.Item(3, 15).Value = 0.5209011059384622
$ws.Cells.Item(3, 16).Value = 0.5209011059384622
$ws.Cells.Item(3, 17).Value = 419.4850960501284
$ws.Cells.Item(3, 18).Value = 3775.365864451156
$ws.Cells.Item(3, 19).Value = 0.1065146398083164
$ws.Cells.Item(3, 20).Value = 0.1065146398083164

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 10.32689766666667
$ws.Cells.Item(4, 8).Value = 30.980693
$ws.Cells.Item(4, 9).Value = 0.2044815006034941
$ws.Cells.Item(4, 10).Value = 0.204481500603494
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 20.825229
$ws.Cells.Item(4, 14).Value = 62.475687
$ws.Cells.Item(4, 15).Value = 0.2670535794123827
$ws.Cells.Item(4, 16).Value = 0.2670535794123827
$ws.Cells.Item(4, 17).Value = 215.060008767899
$ws.Cells.Item(4, 18).Value = 1935.540078911091
$ws.Cells.Item(4, 19).Value = 0.05460751665977838
$ws.Cells.Item(4, 20).Value = 0.05460751665977837

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 29.37031933333333
$ws.Cells.Item(5, 8).Value = 88.110958
$ws.Cells.Item(5, 9).Value = 0.5815577111671272
$ws.Cells.Item(5, 10).Value = 0.5815577111671272
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 16.535604
$ws.Cells.Item(5, 14).Value = 49.606812
$ws.Cells.Item(5, 15).Value = 0.2120453146491552
$ws.Cells.Item(5, 16).Value = 0.2120453146491552
$ws.Cells.Item(5, 17).Value = 485.655969849544
$ws.Cells.Item(5, 18).Value = 4370.903728645896
$ws.Cells.Item(5, 19).Value = 0.123316587851076
$ws.Cells.Item(5, 20).Value = 0.123316587851076

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 29.37031933333333
$ws.Cells.Item(6, 8).Value = 88.110958
$ws.Cells.Item(6, 9).Value = 0.5815577111671272
$ws.Cells.Item(6, 10).Value = 0.5815577111671272
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 40.62063066666667
$ws.Cells.Item(6, 14).Value = 121.861892
$ws.Cells.Item(6, 15).Value = 0.5209011059384622
$ws.Cells.Item(6, 16).Value = 0.5209011059384622
$ws.Cells.Item(6, 17).Value = 1193.040894201393
$ws.Cells.Item(6, 18).Value = 10737.36804781254
$ws.Cells.Item(6, 19).Value = 0.3029340549139973
$ws.Cells.Item(6, 20).Value = 0.3029340549139973

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 29.37031933333333
$ws.Cells.Item(7, 8).Value = 88.110958
$ws.Cells.Item(7, 9).Value = 0.5815577111671272
$ws.Cells.Item(7, 10).Value = 0.5815577111671272
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 20.825229
$ws.Cells.Item(7, 14).Value = 62.475687
$ws.Cells.Item(7, 15).Value = 0.2670535794123827
$ws.Cells.Item(7, 16).Value = 0.2670535794123827
$ws.Cells.Item(7, 17).Value = 611.643625919794
$ws.Cells.Item(7, 18).Value = 5504.792633278146
$ws.Cells.Item(7, 19).Value = 0.1553070684020539
$ws.Cells.Item(7, 20).Value = 0.1553070684020539

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 10.80562866666667
$ws.Cells.Item(8, 8).Value = 32.416886
$ws.Cells.Item(8, 9).Value = 0.2139607882293788
$ws.Cells.Item(8, 10).Value = 0.2139607882293788
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 16.535604
$ws.Cells.Item(8, 14).Value = 49.606812
$ws.Cells.Item(8, 15).Value = 0.2120453146491552
$ws.Cells.Item(8, 16).Value = 0.2120453146491552
$ws.Cells.Item(8, 17).Value = 178.677596603048
$ws.Cells.Item(8, 18).Value = 1608.098369427432
$ws.Cells.Item(8, 19).Value = 0.0453693826626799
$ws.Cells.Item(8, 20).Value = 0.0453693826626799

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 10.80562866666667
$ws.Cells.Item(9, 8).Value = 32.416886
$ws.Cells.Item(9, 9).Value = 0.2139607882293788
$ws.Cells.Item(9, 10).Value = 0.2139607882293788
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 40.62063066666667
$ws.Cells.Item(9, 14).Value = 121.861892
$ws.Cells.Item(9, 15).Value = 0.5209011059384622
$ws.Cells.Item(9, 16).Value = 0.5209011059384622
$ws.Cells.Item(9, 17).Value = 438.9314511898124
$ws.Cells.Item(9, 18).Value = 3950.383060708312
$ws.Cells.Item(9, 19).Value = 0.1114524112161485
$ws.Cells.Item(9, 20).Value = 0.1114524112161485

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 10.80562866666667
$ws.Cells.Item(10, 8).Value = 32.416886
$ws.Cells.Item(10, 9).Value = 0.2139607882293788
$ws.Cells.Item(10, 10).Value = 0.2139607882293788
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 20.825229
$ws.Cells.Item(10, 14).Value = 62.475687
$ws.Cells.Item(10, 15).Value = 0.2670535794123827
$ws.Cells.Item(10, 16).Value = 0.2670535794123827
$ws.Cells.Item(10, 17).Value = 225.029691472298
$ws.Cells.Item(10, 18).Value = 2025.267223250682
$ws.Cells.Item(10, 19).Value = 0.05713899435055041
$ws.Cells.Item(10, 20).Value = 0.0571389943505504
